$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Developers Secure `$224.3 Million in Bridge Financing for The Italic at 26-32 Jackson Avenue In Long Island City, Queens"
$ws.Cells.Item(2, 2).Value = "https://newyorkyimby.com/2026/01/developers-secure-224-3-million-in-bridge-financing-for-the-italic-at-26-32-jackson-avenue-in-long-island-city-queens.html"
$ws.Cells.Item(2, 3).Value = 'Developer American Lions has secured $224.3 million in bridge financing for The Italic, a <a href="https://newyorkyimby.com/2025/03/the-italic-completes-construction-at-26-32-jackson-avenue-in-long-island-city-queens.html">recently completed</a> residential skyscraper at 26-32 Jackson Avenue in <a href="https://newyorkyimby.com/neighborhoods/long-island-city">Long Island City</a>, Queens. Designed by SLCE Architects, the 49-story tower yields 363 residential units, including 109 affordable apartments. JLL Capital Markets arranged the financing package for the development team, a joint venture between The Carlyle Group, Fetner Properties, and Lions Group. The property is located at the corner of Jackson Avenue and Purves Street.'
$ws.Cells.Item(2, 4).Value = "2026-01-28T12:30:04+00:00"
$ws.Cells.Item(2, 5).Value = "Wed, 28 Jan 2026 12:30:04 +0000"
$ws.Cells.Item(2, 6).Value = "YIMBY"
$ws.Cells.Item(2, 7).Value = "YIMBY - Long Island City"

# Column H has no content for this row, but the source row still has a
# (text-typed, empty) cell present at H2 so the sheet's dimension grows
# to A1:H2. Touch the cell's formatting so it is materialized even
# though it carries no value, then reset it back to the default style.
$ws.Cells.Item(2, 8).NumberFormat = "0"
$ws.Cells.Item(2, 8).Style = "Normal"

$wb.Save()
